# Update Name of Algo
# Applies updated imputed values (result_data_KNN) to the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -13.182
$ws.Range("E3").Value = 16.228
$ws.Range("A12").Value = -21.651
$ws.Range("C14").Value = -12.139
$ws.Range("C26").Value = -12.854
$ws.Range("E30").Value = 16.135
$ws.Range("C31").Value = -13.186
$ws.Range("A32").Value = -21.742
$ws.Range("C35").Value = -12.762
$ws.Range("A36").Value = -20.178
$ws.Range("C37").Value = -13.341
$ws.Range("A38").Value = -19.741
$ws.Range("E44").Value = 16.611
$ws.Range("C45").Value = -12.898
$ws.Range("A46").Value = -21.92
$ws.Range("A54").Value = -22.155
$ws.Range("A55").Value = -22.145
$ws.Range("C57").Value = -13.829
$ws.Range("E58").Value = 16.536
$ws.Range("A67").Value = -21.495
$ws.Range("A69").Value = -21.651
$ws.Range("A72").Value = -21.445
$ws.Range("E84").Value = 16.359
$ws.Range("E89").Value = 16.945
$ws.Range("A91").Value = -21.508
$ws.Range("E91").Value = 17.39
$ws.Range("E92").Value = 17.381
$ws.Range("A99").Value = -20.428
$ws.Range("C100").Value = -12.578
$ws.Range("C102").Value = -12.87
$ws.Range("E102").Value = 16.418
